$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update "Bad Drivers" summary row (row 3) and Totals row (row 4) ---
$ws.Range("C3").Value = 450
$ws.Range("D3").Value = 91.3
$ws.Range("C4").Value = 450

# --- 2) Insert a new "Good Driver" row before the existing row 12 ---
# (everything currently at rows 12-17 shifts down to rows 13-18)
$ws.Rows.Item(12).Insert()

# --- 3) Populate the newly inserted row 12 with the new driver entry ---
$ws.Range("A12").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Range("B12").NumberFormat = "#,##0"
$ws.Range("B12").Value = 11128
$ws.Range("D12").Value = 100
$ws.Range("E12").Value = 0

# --- 4) Refresh "Total Samples" (column B) for the rows that shifted down ---
$ws.Range("B13").Value = 486214
$ws.Range("B14").Value = 79953
$ws.Range("B15").Value = 35355
$ws.Range("B16").Value = 65425
$ws.Range("B17").Value = 117653
$ws.Range("B18").Value = 56018
